$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.132470666666666
$ws.Range("H2").Value = 21.397412
$ws.Range("I2").Value = 0.1078130252899183
$ws.Range("J2").Value = 0.1078130252899183
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 274.9758354464844
$ws.Range("R2").Value = 2474.78251901836
$ws.Range("S2").Value = 0.06208037576454386
$ws.Range("T2").Value = 0.06208037576454386

# Row 3
$ws.Range("G3").Value = 7.132470666666666
$ws.Range("H3").Value = 21.397412
$ws.Range("I3").Value = 0.1078130252899183
$ws.Range("J3").Value = 0.1078130252899183
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 38.30618903017067
$ws.Range("R3").Value = 344.7557012715361
$ws.Range("S3").Value = 0.008648260328909736
$ws.Range("T3").Value = 0.008648260328909738

# Row 4
$ws.Range("G4").Value = 7.132470666666666
$ws.Range("H4").Value = 21.397412
$ws.Range("I4").Value = 0.1078130252899183
$ws.Range("J4").Value = 0.1078130252899183
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 164.2598127948911
$ws.Range("R4").Value = 1478.33831515402
$ws.Range("S4").Value = 0.03708438919646472
$ws.Range("T4").Value = 0.03708438919646473

# Row 5
$ws.Range("I5").Value = 0.2490596131114117
$ws.Range("J5").Value = 0.2490596131114118
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 635.2235734701411
$ws.Range("R5").Value = 5717.012161231271
$ws.Range("S5").Value = 0.1434123041084368
$ws.Range("T5").Value = 0.1434123041084368

# Row 6
$ws.Range("I6").Value = 0.2490596131114117
$ws.Range("J6").Value = 0.2490596131114118
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 88.4913914062947
$ws.Range("R6").Value = 796.4225226566523
$ws.Range("S6").Value = 0.01997840581704226
$ws.Range("T6").Value = 0.01997840581704227

# Row 7
$ws.Range("I7").Value = 0.2490596131114117
$ws.Range("J7").Value = 0.2490596131114118
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 379.4577261368628
$ws.Range("R7").Value = 3415.119535231765
$ws.Range("S7").Value = 0.08566890318593263
$ws.Range("T7").Value = 0.08566890318593266

# Row 8
$ws.Range("G8").Value = 42.546687
$ws.Range("H8").Value = 127.640061
$ws.Range("I8").Value = 0.6431273615986699
$ws.Range("J8").Value = 0.6431273615986699
$ws.Range("M8").Value = 38.55267666666666
$ws.Range("N8").Value = 115.65803
$ws.Range("O8").Value = 0.5758151725879548
$ws.Range("P8").Value = 0.5758151725879548
$ws.Range("Q8").Value = 1640.28866714887
$ws.Range("R8").Value = 14762.59800433983
$ws.Range("S8").Value = 0.3703224927149741
$ws.Range("T8").Value = 0.3703224927149741

# Row 9
$ws.Range("G9").Value = 42.546687
$ws.Range("H9").Value = 127.640061
$ws.Range("I9").Value = 0.6431273615986699
$ws.Range("J9").Value = 0.6431273615986699
$ws.Range("O9").Value = 0.08021535714867321
$ws.Range("P9").Value = 0.08021535714867323
$ws.Range("Q9").Value = 228.504470750412
$ws.Range("R9").Value = 2056.540236753708
$ws.Range("S9").Value = 0.05158869100272121
$ws.Range("T9").Value = 0.05158869100272122

# Row 10
$ws.Range("G10").Value = 42.546687
$ws.Range("H10").Value = 127.640061
$ws.Range("I10").Value = 0.6431273615986699
$ws.Range("J10").Value = 0.6431273615986699
$ws.Range("M10").Value = 23.02986166666667
$ws.Range("N10").Value = 69.089585
$ws.Range("O10").Value = 0.3439694702633719
$ws.Range("P10").Value = 0.3439694702633719
$ws.Range("Q10").Value = 979.844315984965
$ws.Range("R10").Value = 8818.598843864684
$ws.Range("S10").Value = 0.2212161778809745
$ws.Range("T10").Value = 0.2212161778809745
